# DNB Mastercard Demo -> generic "Sheet" with refreshed / expanded transaction data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab (was "DNB Mastercard Demo").
$ws.Name = "Sheet"

# 2. New transaction rows (dates are Excel serial numbers, amounts are "Ut"
#    unless noted). Row 18 ("Innbetaling") posts to "Inn" (column E) instead.
#    This grows the sheet from 14 to 20 data rows and swaps in a more
#    diverse / realistic set of merchants.
$rows = @(
    @{ Row = 2;  Date = 45747; Desc = "SPOTIFY";                  Out = 129 },
    @{ Row = 3;  Date = 45746; Desc = "Kiwi Torshov";              Out = 312 },
    @{ Row = 4;  Date = 45744; Desc = "netflix.com";               Out = 179 },
    @{ Row = 5;  Date = 45743; Desc = "kiwi grünerløkka";          Out = 267.5 },
    @{ Row = 6;  Date = 45741; Desc = "NILLE STORO";                Out = 149 },
    @{ Row = 7;  Date = 45740; Desc = "VINMONOPOLET GRÜNERLØKKA";   Out = 567 },
    @{ Row = 8;  Date = 45738; Desc = "Starbucks Aker Brygge";      Out = 89 },
    @{ Row = 9;  Date = 45736; Desc = "NORMAL MAJORSTUEN";          Out = 199 },
    @{ Row = 10; Date = 45734; Desc = "GITHUB.COM";                 Out = 129 },
    @{ Row = 11; Date = 45733; Desc = "FLYING TIGER OSLO";          Out = 89 },
    @{ Row = 12; Date = 45731; Desc = "REMA 1000 GRÜNERLØKKA";      Out = 534.2 },
    @{ Row = 13; Date = 45730; Desc = "starbucks bogstadveien";     Out = 75 },
    @{ Row = 14; Date = 45728; Desc = "MENY MAJORSTUEN";            Out = 623.45 },
    @{ Row = 15; Date = 45726; Desc = "KICKS OSLO CITY";            Out = 456 },
    @{ Row = 16; Date = 45724; Desc = "INTERSPORT CC VEST";         Out = 1299 },
    @{ Row = 17; Date = 45722; Desc = "VITA KARL JOHAN";            Out = 289 },
    @{ Row = 18; Date = 45721; Desc = "Innbetaling";                In  = 15000 },
    @{ Row = 19; Date = 45719; Desc = "COOP MEGA TORSHOV";          Out = 756.4 },
    @{ Row = 20; Date = 45717; Desc = "ESSO HOVINBYEN";             Out = 678 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Date
    # 3. Widen the date format to include a time component
    #    ("yyyy-mm-dd" -> "yyyy-mm-dd h:mm:ss") on every date cell.
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
    $ws.Cells.Item($row, 2).Value = $r.Desc

    # Clear stale E/F values left over from the old template before writing
    # the new ones, since rows are being re-purposed.
    $ws.Cells.Item($row, 5).Value = $null
    $ws.Cells.Item($row, 6).Value = $null

    if ($r.ContainsKey("In")) {
        $ws.Cells.Item($row, 5).Value = $r.In
    }
    if ($r.ContainsKey("Out")) {
        $ws.Cells.Item($row, 6).Value = $r.Out
    }
}
